$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to keep text formatting so numeric-looking
# strings (e.g. "1.003") are not auto-converted to floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.371.90'
$ws.Range("E2").Value = '  -1.56%  '

$ws.Range("D3").Value = '1.828.02'
$ws.Range("E3").Value = '  -1.65%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -3.31%  '

$ws.Range("D5").Value = '315.56'
$ws.Range("E5").Value = '  -2.45%  '

$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -2.90%  '

$ws.Range("D7").Value = '0.4299'
$ws.Range("E7").Value = '  -2.55%  '

$ws.Range("D8").Value = '0.3701'
$ws.Range("E8").Value = '  -3.08%  '

$ws.Range("D9").Value = '0.07257'
$ws.Range("E9").Value = '  -2.55%  '

$ws.Range("D10").Value = '0.8661'
$ws.Range("E10").Value = '  -2.55%  '

$ws.Range("D11").Value = '21.16'
$ws.Range("E11").Value = '  -2.07%  '

$ws.Range("D12").Value = '1.824.63'
$ws.Range("E12").Value = '  -1.87%  '

$ws.Range("D13").Value = '6.680'
$ws.Range("E13").Value = '  -1.03%  '

$ws.Range("D14").Value = '5.356'
$ws.Range("E14").Value = '  -3.33%  '

$ws.Range("D15").Value = '0.07104'
$ws.Range("E15").Value = '  -1.51%  '

$ws.Range("D16").Value = '87.86'
$ws.Range("E16").Value = '  +2.01%  '

$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  -3.31%  '

$ws.Range("D18").Value = '0.000008897'
$ws.Range("E18").Value = '  -2.43%  '

$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -3.05%  '

$ws.Range("D20").Value = '15.23'
$ws.Range("E20").Value = '  -2.24%  '

$ws.Range("D21").Value = '27.380.96'
$ws.Range("E21").Value = '  -1.56%  '

$ws.Range("D22").Value = '5.162'
$ws.Range("E22").Value = '  -2.75%  '

$ws.Range("D23").Value = '10.87'
$ws.Range("E23").Value = '  -3.63%  '

$ws.Range("D24").Value = '2.051.01'
$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").Value = '2.011'
$ws.Range("E25").Value = '  -3.08%  '

$ws.Range("D26").Value = '153.27'
$ws.Range("E26").Value = '  -3.74%  '

$ws.Range("D27").Value = '18.47'
$ws.Range("E27").Value = '  -1.67%  '

$ws.Range("D28").Value = '2.140'
$ws.Range("E28").Value = '  +7.17%  '

$ws.Range("D29").Value = '5.298'
$ws.Range("E29").Value = '  -1.00%  '

$ws.Range("D30").Value = '116.98'
$ws.Range("E30").Value = '  -1.64%  '

$ws.Range("D31").Value = '0.08857'
$ws.Range("E31").Value = '  -2.91%  '

$ws.Range("D32").Value = '1.206'
$ws.Range("E32").Value = '  -0.87%  '

$ws.Range("D33").Value = '0.7665'
$ws.Range("E33").Value = '  -0.89%  '

$ws.Range("D34").Value = '4.507'
$ws.Range("E34").Value = '  -2.24%  '

$ws.Range("D35").Value = '2.847'
$ws.Range("E35").Value = '  -5.98%  '

$ws.Range("D36").Value = '1.003'
$ws.Range("E36").Value = '  -3.12%  '

$ws.Range("D37").Value = '1.121'
$ws.Range("E37").Value = '  -3.18%  '

$ws.Range("D38").Value = '0.01960'
$ws.Range("E38").Value = '  -1.38%  '

$ws.Range("E39").Value = '  -0.92%  '

$ws.Range("D40").Value = '2.881'
$ws.Range("E40").Value = '  +0.73%  '

$ws.Range("D41").Value = '7.129'
$ws.Range("E41").Value = '  +2.41%  '

$ws.Range("D42").Value = '0.1682'
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("E43").Value = '  -2.88%  '

$ws.Range("D44").Value = '8.671'
$ws.Range("E44").Value = '  -1.54%  '

$ws.Range("E45").Value = '  -2.75%  '

$ws.Range("D46").Value = '106.34'
$ws.Range("E46").Value = '  -4.09%  '

$ws.Range("D47").Value = '0.4731'
$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").Value = '0.06421'
$ws.Range("E48").Value = '  -2.40%  '

$ws.Range("D49").Value = '1.004'
$ws.Range("E49").Value = '  -3.19%  '

$ws.Range("D50").Value = '1.674'
$ws.Range("E50").Value = '  -2.49%  '

$ws.Range("D51").Value = '1.825'
$ws.Range("E51").Value = '  -3.27%  '

# Restore the default (Normal) cell style now that the text values are set,
# matching the original workbook which had no explicit style on these cells.
$ws.Range("D2:E51").Style = "Normal"

Write-Output "Updated cryptos list values"